$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply updated odds values per the diff (row-by-row, column order)

# Row 4
$ws.Range("U4").Value = 1.65
$ws.Range("V4").Value = 2.19
$ws.Range("W4").Value = 6.4
$ws.Range("X4").Value = 7
$ws.Range("AD4").Value = 6
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 20

# Row 7
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 3.45
$ws.Range("I7").Value = 5.4
$ws.Range("P7").Value = 2.8
$ws.Range("T7").Value = 2.57
$ws.Range("U7").Value = 1.93
$ws.Range("X7").Value = 6.9
$ws.Range("Z7").Value = 12
$ws.Range("AB7").Value = 32
$ws.Range("AD7").Value = 6.9
$ws.Range("AE7").Value = 18
$ws.Range("AF7").Value = 100
$ws.Range("AN7").Value = 3.35
$ws.Range("AP7").Value = 18
$ws.Range("AR7").Value = 60
$ws.Range("AU7").Value = 7.6
$ws.Range("AY7").Value = 35

# Row 8
$ws.Range("G8").Value = 1.55
$ws.Range("H8").Value = 3.85
$ws.Range("I8").Value = 5.4
$ws.Range("J8").Value = 2.07
$ws.Range("K8").Value = 2.27
$ws.Range("L8").Value = 5.1
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 12.7
$ws.Range("P8").Value = 3.65
$ws.Range("Q8").Value = 1.55
$ws.Range("R8").Value = 2.15
$ws.Range("S8").Value = 1.29
$ws.Range("T8").Value = 3.32
$ws.Range("W8").Value = 8.25
$ws.Range("X8").Value = 8.25
$ws.Range("Y8").Value = 7.9
$ws.Range("Z8").Value = 12
$ws.Range("AA8").Value = 11.5
$ws.Range("AB8").Value = 20
$ws.Range("AC8").Value = 14
$ws.Range("AD8").Value = 7.9
$ws.Range("AE8").Value = 13.5
$ws.Range("AH8").Value = 19.5
$ws.Range("AI8").Value = 40
$ws.Range("AJ8").Value = 16.5
$ws.Range("AK8").Value = 110
$ws.Range("AL8").Value = 50
$ws.Range("AM8").Value = 40
$ws.Range("AN8").Value = 3.55
$ws.Range("AO8").Value = 7.4
$ws.Range("AP8").Value = 14.5
$ws.Range("AQ8").Value = 22
$ws.Range("AR8").Value = 45
$ws.Range("AS8").Value = 150
$ws.Range("AT8").Value = 3.1
$ws.Range("AU8").Value = 6.8
$ws.Range("AW8").Value = 7.1
$ws.Range("AX8").Value = 28
$ws.Range("AY8").Value = 27
$ws.Range("BB8").Value = 250

# Row 9
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("Q9").Value = 2.08
$ws.Range("R9").Value = 1.73
